# Updated symbol list (Price/Volume(1h) columns) to refresh the crypto
# price snapshot. Price (D) and Volume 1h % (E) columns are stored as
# plain text in this sheet, so force Text number-format before writing
# the new values to keep them as text rather than have Excel coerce
# them into numbers/percentages.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:E2").NumberFormat = "@"
$ws.Range("D2").Value = "329.22"
$ws.Range("E2").Value = "6.29%"

$ws.Range("D3:E3").NumberFormat = "@"
$ws.Range("D3").Value = "40.51"
$ws.Range("E3").Value = "10.44%"

$ws.Range("D4:E4").NumberFormat = "@"
$ws.Range("D4").Value = "6.006"
$ws.Range("E4").Value = "17.78%"

$ws.Range("D5:E5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08144"
$ws.Range("E5").Value = "5.66%"

$ws.Range("D6:E6").NumberFormat = "@"
$ws.Range("D6").Value = "4.611"
$ws.Range("E6").Value = "5.00%"

$ws.Range("D7:E7").NumberFormat = "@"
$ws.Range("D7").Value = "8.769"
$ws.Range("E7").Value = "5.33%"

$ws.Range("D8:E8").NumberFormat = "@"
$ws.Range("D8").Value = "1.968"
$ws.Range("E8").Value = "6.26%"

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-0.20%"

$ws.Range("D10:E10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9481"
$ws.Range("E10").Value = "2.59%"

$ws.Range("D11:E11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1330"
$ws.Range("E11").Value = "15.46%"

$ws.Range("D12:E12").NumberFormat = "@"
$ws.Range("D12").Value = "0.1995"
$ws.Range("E12").Value = "6.24%"

$ws.Range("D13:E13").NumberFormat = "@"
$ws.Range("D13").Value = "9.693"
$ws.Range("E13").Value = "53.85%"

$ws.Range("D14:E14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09296"
$ws.Range("E14").Value = "5.57%"

$ws.Range("D15:E15").NumberFormat = "@"
$ws.Range("D15").Value = "0.03496"
$ws.Range("E15").Value = "4.02%"

$ws.Range("D16:E16").NumberFormat = "@"
$ws.Range("D16").Value = "0.09612"
$ws.Range("E16").Value = "0.83%"

$ws.Range("D17:E17").NumberFormat = "@"
$ws.Range("D17").Value = "0.001333"
$ws.Range("E17").Value = "-3.58%"

$ws.Range("D18:E18").NumberFormat = "@"
$ws.Range("D18").Value = "0.006430"
$ws.Range("E18").Value = "7.52%"

$ws.Range("D19:E19").NumberFormat = "@"
$ws.Range("D19").Value = "3.355"
$ws.Range("E19").Value = "-0.10%"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.3543"

$ws.Range("D21:E21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1414"
$ws.Range("E21").Value = "9.52%"

$ws.Range("D22:E22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2415"
$ws.Range("E22").Value = "4.52%"

$ws.Range("D23:E23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04439"
$ws.Range("E23").Value = "2.20%"

$ws.Range("D24:E24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001259"
$ws.Range("E24").Value = "4.67%"

$ws.Range("D25:E25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004438"
$ws.Range("E25").Value = "4.43%"

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-17.85%"

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "4.65%"

$ws.Range("D39:E39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02460"
$ws.Range("E39").Value = "15.88%"

$ws.Range("D40:E40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05293"
$ws.Range("E40").Value = "5.39%"

$ws.Range("D41:E41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007440"
$ws.Range("E41").Value = "-0.94%"

$ws.Range("D42:E42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1436"
$ws.Range("E42").Value = "6.48%"

$ws.Range("D43:E43").NumberFormat = "@"
$ws.Range("D43").Value = "0.009009"
$ws.Range("E43").Value = "7.02%"

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-0.68%"

$ws.Range("D45:E45").NumberFormat = "@"
$ws.Range("D45").Value = "0.01060"
$ws.Range("E45").Value = "37.23%"

$ws.Range("D46:E46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006824"
$ws.Range("E46").Value = "7.29%"

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.18%"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003504"

$ws.Range("D49:E49").NumberFormat = "@"
$ws.Range("D49").Value = "0.001805"
$ws.Range("E49").Value = "6.74%"

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.18%"

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.18%"
